# save data done + era data updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) — header styled like the other header cells (G1 etc.)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Per-row save flags
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
